$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values -- must stay text even when numeric-looking,
# matching the original inlineStr cell type. Force text via NumberFormat,
# write the value, then restore the cell's style so no stray numFmt/style
# index is left behind.
$priceUpdates = @{
    "D2"  = "42.314.99"
    "D3"  = "2.296.08"
    "D5"  = "315.83"
    "D6"  = "104.35"
    "D10" = "39.67"
    "D12" = "8.30"
    "D14" = "0.963"
    "D15" = "15.31"
    "D16" = "2.644.91"
    "D17" = "2.298.75"
    "D18" = "42.269.75"
    "D19" = "7.38"
    "D21" = "73.44"
    "D23" = "277.73"
    "D24" = "10.46"
    "D27" = "10.85"
    "D28" = "2.36"
    "D29" = "22.84"
    "D30" = "36.32"
    "D31" = "164.23"
    "D33" = "5.83"
    "D34" = "2.82"
    "D38" = "0.0349"
    "D39" = "3.77"
    "D40" = "2.80"
    "D41" = "99.38"
    "D43" = "69.45"
    "D46" = "12.02"
    "D47" = "112.04"
    "D48" = "77.54"
    "D49" = "8.93"
    "D51" = "1.589.10"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Column E ("Volume(1h)") values -- plain percentage-with-padding strings.
$volumeUpdates = @{
    "E2"  = "  -1.48%  "
    "E3"  = "  -3.02%  "
    "E4"  = "  -0.02%  "
    "E5"  = "  -0.44%  "
    "E6"  = "  -4.41%  "
    "E8"  = "  +0.00%  "
    "E9"  = "  -2.75%  "
    "E10" = "  -4.60%  "
    "E11" = "  -2.54%  "
    "E12" = "  -3.97%  "
    "E13" = "  -0.36%  "
    "E14" = "  -4.96%  "
    "E15" = "  -4.73%  "
    "E16" = "  -2.91%  "
    "E17" = "  -2.99%  "
    "E18" = "  -1.52%  "
    "E19" = "  -3.88%  "
    "E20" = "  -1.15%  "
    "E21" = "  -3.96%  "
    "E22" = "  -1.00%  "
    "E23" = "  +3.78%  "
    "E24" = "  +9.81%  "
    "E25" = "  -2.39%  "
    "E26" = "  +0.60%  "
    "E27" = "  -5.72%  "
    "E28" = "  +4.95%  "
    "E29" = "  -2.68%  "
    "E30" = "  -2.76%  "
    "E31" = "  -2.90%  "
    "E32" = "  -4.16%  "
    "E33" = "  -3.86%  "
    "E34" = "  -2.94%  "
    "E37" = "  -4.02%  "
    "E38" = "  -4.33%  "
    "E39" = "  -3.21%  "
    "E40" = "  +3.58%  "
    "E41" = "  -7.00%  "
    "E42" = "  -4.39%  "
    "E43" = "  -2.96%  "
    "E44" = "  -5.69%  "
    "E45" = "  +0.11%  "
    "E46" = "  -3.97%  "
    "E47" = "  -2.42%  "
    "E48" = "  -3.98%  "
    "E49" = "  -2.89%  "
    "E50" = "  -5.16%  "
    "E51" = "  +0.17%  "
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
